$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the data columns so numeric-looking strings
# (e.g. "1.010", "26.906.47") are not auto-converted to numbers by Excel,
# matching the inline-string (text) cells in the source workbook.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.906.47'
$ws.Range("E2").Value = '  -2.21%  '

# Row 3
$ws.Range("D3").Value = '1.817.45'
$ws.Range("E3").Value = '  -1.45%  '

# Row 4
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  -0.47%  '

# Row 5
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.009'
$ws.Range("E5").Value = '  -0.26%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '308.82'
$ws.Range("E6").Value = '  -2.04%  '

# Row 7
$ws.Range("D7").Value = '0.4580'
$ws.Range("E7").Value = '  -3.54%  '

# Row 8
$ws.Range("D8").Value = '0.3625'
$ws.Range("E8").Value = '  -2.30%  '

# Row 9
$ws.Range("D9").Value = '0.07192'
$ws.Range("E9").Value = '  -3.86%  '

# Row 10
$ws.Range("D10").Value = '0.8559'
$ws.Range("E10").Value = '  -3.62%  '

# Row 11
$ws.Range("D11").Value = '19.68'
$ws.Range("E11").Value = '  -4.19%  '

# Row 12
$ws.Range("D12").Value = '1.877.14'
$ws.Range("E12").Value = '  +1.19%  '

# Row 13
$ws.Range("D13").Value = '0.07564'
$ws.Range("E13").Value = '  +2.42%  '

# Row 14
$ws.Range("D14").Value = '92.70'
$ws.Range("E14").Value = '  -0.99%  '

# Row 15
$ws.Range("D15").Value = '5.310'
$ws.Range("E15").Value = '  -3.15%  '

# Row 16
$ws.Range("D16").Value = '6.475'
$ws.Range("E16").Value = '  -2.00%  '

# Row 17
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  -0.36%  '

# Row 18
$ws.Range("D18").Value = '0.000008594'
$ws.Range("E18").Value = '  -3.05%  '

# Row 19
$ws.Range("D19").Value = '1.010'
$ws.Range("E19").Value = '  -0.12%  '

# Row 20
$ws.Range("D20").Value = '27.215.17'
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("D21").Value = '14.38'
$ws.Range("E21").Value = '  -3.23%  '

# Row 22
$ws.Range("D22").Value = '5.121'
$ws.Range("E22").Value = '  -4.14%  '

# Row 23
$ws.Range("D23").Value = '10.47'
$ws.Range("E23").Value = '  -2.36%  '

# Row 24
$ws.Range("D24").Value = '2.100.07'
$ws.Range("E24").Value = '  +0.93%  '

# Row 25
$ws.Range("D25").Value = '151.20'
$ws.Range("E25").Value = '  -0.66%  '

# Row 26
$ws.Range("D26").Value = '1.856'
$ws.Range("E26").Value = '  -2.67%  '

# Row 27
$ws.Range("D27").Value = '18.10'
$ws.Range("E27").Value = '  -2.96%  '

# Row 28
$ws.Range("D28").Value = '2.073'
$ws.Range("E28").Value = '  -4.68%  '

# Row 29
$ws.Range("D29").Value = '5.056'
$ws.Range("E29").Value = '  -4.29%  '

# Row 30
$ws.Range("D30").Value = '115.27'
$ws.Range("E30").Value = '  -2.62%  '

# Row 31
$ws.Range("D31").Value = '0.08846'
$ws.Range("E31").Value = '  -1.83%  '

# Row 32
$ws.Range("D32").Value = '2.950'
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").Value = '1.133'
$ws.Range("E33").Value = '  -4.30%  '

# Row 34
$ws.Range("D34").Value = '0.7179'
$ws.Range("E34").Value = '  -5.66%  '

# Row 35
$ws.Range("D35").Value = '4.382'
$ws.Range("E35").Value = '  -4.17%  '

# Row 36
$ws.Range("D36").Value = '1.010'
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("D37").Value = '2.461'
$ws.Range("E37").Value = '  +2.23%  '

# Row 38
$ws.Range("D38").Value = '1.072'
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("D39").Value = '0.05230'
$ws.Range("E39").Value = '  -2.29%  '

# Row 40
$ws.Range("D40").Value = '0.01902'
$ws.Range("E40").Value = '  -3.27%  '

# Row 41
$ws.Range("D41").Value = '2.913'
$ws.Range("E41").Value = '  -2.98%  '

# Row 42
$ws.Range("D42").Value = '7.098'
$ws.Range("E42").Value = '  -3.54%  '

# Row 43
$ws.Range("D43").Value = '0.5111'
$ws.Range("E43").Value = '  -4.62%  '

# Row 44
$ws.Range("D44").Value = '0.1618'
$ws.Range("E44").Value = '  -2.80%  '

# Row 45
$ws.Range("D45").Value = '8.172'
$ws.Range("E45").Value = '  -4.66%  '

# Row 46
$ws.Range("D46").Value = '0.4776'
$ws.Range("E46").Value = '  -3.78%  '

# Row 47
$ws.Range("D47").Value = '1.010'
$ws.Range("E47").Value = '  -0.27%  '

# Row 48
$ws.Range("E48").Value = '  -4.48%  '

# Row 49
$ws.Range("D49").Value = '102.60'
$ws.Range("E49").Value = '  -2.37%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06218'
$ws.Range("E50").Value = '  -1.75%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.606'
$ws.Range("E51").Value = '  -4.61%  '
